$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

# Schedule counter update: bump the auto-test counter after a run
$ws.Range("A2").Value = "AGSAutoT03"

# Move the selection to A2
$ws.Range("A2").Select()

# Column B should match column A's width (not auto bestFit anymore)
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth
